$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Top language-picker line (paragraph 1): "English" (hyperlink run)
#    and the " / Portuguese / French / ..." run.
# -----------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1.Range.Find.Execute("English", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Anglais", 2)
$p1.Range.Find.Execute(" / Portuguese / French / Thai / Vietnamese / Spanish", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " / portugais / français / thaïlandais / vietnamien / espagnol", 2)

# -----------------------------------------------------------------
# 2) Second "English" heading (paragraph 3).
# -----------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$p3.Range.Find.Execute("English", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Anglais", 2)

# -----------------------------------------------------------------
# 3) Table cell - "Brief:" (paragraph 5). Keep "Brief"/":" as two
#    separate runs by replacing each inside the paragraph scope only.
# -----------------------------------------------------------------
$p5 = $d.Paragraphs(5)
$p5.Range.Find.Execute("Brief", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Résumé", 2)
$p5.Range.Find.Execute(":", $true, $false, $false, $false, $false, `
    $true, 1, $false, " :", 2)

# -----------------------------------------------------------------
# 4) Table cell - brief description (paragraph 6).
# -----------------------------------------------------------------
$p6 = $d.Paragraphs(6)
$p6.Range.Find.Execute("It will be sent via customer.io", $true, $false, `
    $false, $false, $false, $true, 1, $false, "Il sera envoyé via customer.io", 2)

# -----------------------------------------------------------------
# 5) Table cell - "Target audience:" (paragraph 8). Same two-run
#    handling as "Brief:".
# -----------------------------------------------------------------
$p8 = $d.Paragraphs(8)
$p8.Range.Find.Execute("Target audience", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Public cible", 2)
$p8.Range.Find.Execute(":", $true, $false, $false, $false, $false, `
    $true, 1, $false, " :", 2)

# -----------------------------------------------------------------
# 6) "We'll miss you at the [EVENT NAME]!" heading (paragraph 14).
#    Only this [EVENT NAME] (and the one in paragraph 18) change -
#    the one in the Subject line paragraph stays in English.
# -----------------------------------------------------------------
$p14 = $d.Paragraphs(14)
$p14.Range.Find.Execute("We" + [char]0x2019 + "ll miss you at the ", $true, `
    $false, $false, $false, $false, $true, 1, $false, "Vous nous manquerez au ", 2)
$p14.Range.Find.Execute("[EVENT NAME]", $true, $false, $false, $false, $false, `
    $true, 1, $false, "[NOM DE L'ÉVÉNEMENT]", 2)

# -----------------------------------------------------------------
# 7) "Dear [PARTNER NAME], " (paragraph 16).
# -----------------------------------------------------------------
$p16 = $d.Paragraphs(16)
$p16.Range.Find.Execute("Dear ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Cher ", 2)

# -----------------------------------------------------------------
# 8) Thank-you paragraph with the second [EVENT NAME] (paragraph 18).
# -----------------------------------------------------------------
$p18 = $d.Paragraphs(18)
$p18.Range.Find.Execute("Thank you for taking the time to respond to our invitation to the upcoming ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Nous vous remercions d'avoir pris le temps de répondre à notre invitation au prochain ", 2)
$p18.Range.Find.Execute("[EVENT NAME]", $true, $false, $false, $false, $false, `
    $true, 1, $false, "[NOM DE L'ÉVÉNEMENT]", 2)
$p18.Range.Find.Execute(". We were really looking forward to seeing you there.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    ". Nous étions vraiment impatients de vous y retrouver.", 2)

# -----------------------------------------------------------------
# 9) Remaining body paragraphs (19-23).
# -----------------------------------------------------------------
$p19 = $d.Paragraphs(19)
$p19.Range.Find.Execute("Even though we" + [char]0x2019 + "re disappointed we can" + `
    [char]0x2019 + "t meet you, we understand that scheduling conflicts and other commitments sometimes come up. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Bien que nous soyons déçus de ne pas pouvoir vous rencontrer, nous comprenons que des contraintes liées à votre emploi de temps et d'autres engagements peuvent survenir. ", 2)

$p20 = $d.Paragraphs(20)
$p20.Range.Find.Execute("If you" + [char]0x2019 + "re comfortable sharing it with us, we" + `
    [char]0x2019 + "d like to know why you responded no. Please reply to this email as your feedback could help us make improvements in our event planning processes and better serve you in the future.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Si vous souhaitez les partager, nous aimerions connaitre les raisons pour lesquelles vous n'avez pas pu assister à l'événement. Veuillez répondre à ce courriel car vos commentaires pourraient nous aider à améliorer nos processus de planification d'événements et à mieux vous servir à l'avenir.", 2)

$p21 = $d.Paragraphs(21)
$p21.Range.Find.Execute("We hope to see you at our future events. ", $true, `
    $false, $false, $false, $false, $true, 1, $false, `
    "Nous espérons vous rencontrer lors de nos prochains événements. ", 2)

# -----------------------------------------------------------------
# 10) "live chat" / "WhatsApp" paragraph (paragraph 22).
# -----------------------------------------------------------------
$p22 = $d.Paragraphs(22)
$p22.Range.Find.Execute("If you have any questions, please contact us via ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Si vous avez des questions, veuillez nous contacter par ", 2)
$p22.Range.Find.Execute("live chat", $true, $false, $false, $false, $false, `
    $true, 1, $false, "chat en direct", 2)
$p22.Range.Find.Execute(" or ", $true, $false, $false, $false, $false, `
    $true, 1, $false, " ou sur ", 2)

# -----------------------------------------------------------------
# 11) Country-manager paragraph (paragraph 23).
# -----------------------------------------------------------------
$p23 = $d.Paragraphs(23)
$p23.Range.Find.Execute("If you have any questions, please contact your country manager, ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Si vous avez des questions, veuillez contacter votre responsable national, ", 2)
$p23.Range.Find.Execute(", at ", $true, $false, $false, $false, $false, `
    $true, 1, $false, ", à l'adresse ", 2)
$p23.Range.Find.Execute(" or ", $true, $false, $false, $false, $false, `
    $true, 1, $false, " ou au", 2)

# -----------------------------------------------------------------
# 12) Comment text: "choose either one".
# -----------------------------------------------------------------
$c0 = $d.Comments(1)
$c0.Range.Text = "choisissez l'un ou l'autre"
